# Auto-generated Excel COM-interop script
# Applies updated market-data values (columns H-N) for specific Leve rows
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR worksheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 11 (ALC)
$ws.Range("H11").Value = 32.75
$ws.Range("I11").Value = 32.75
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 32.75
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 107.25

# Row 33 (ALC)
$ws.Range("H33").Value = 16868.25
$ws.Range("I33").Value = 247.77777
$ws.Range("J33").Value = 66729.664
$ws.Range("K33").Value = 247.77777
$ws.Range("L33").Value = 66729.664
$ws.Range("M33").Value = -18.77777
$ws.Range("N33").Value = -67187.664

# Row 40 (ALC)
$ws.Range("H40").Value = 5749.5
$ws.Range("I40").Value = 4999
$ws.Range("J40").Value = 5999.6665
$ws.Range("K40").Value = 4999
$ws.Range("L40").Value = 5999.6665
$ws.Range("M40").Value = -4824
$ws.Range("N40").Value = -6349.6665

# Row 86 (ALC)
$ws.Range("H86").Value = 6262
$ws.Range("I86").Value = 6262
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 6262
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -5139

# Row 89 (ALC)
$ws.Range("H89").Value = 6262
$ws.Range("I89").Value = 6262
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 31310
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -25694

# Row 92 (ALC)
$ws.Range("H92").Value = 984.3570999999999
$ws.Range("I92").Value = 984.6667
$ws.Range("J92").Value = 982.5
$ws.Range("K92").Value = 984.6667
$ws.Range("L92").Value = 982.5
$ws.Range("M92").Value = 263.3333
$ws.Range("N92").Value = -3478.5

# Row 103 (ALC)
$ws.Range("H103").Value = 4791.6665
$ws.Range("I103").Value = 5000
$ws.Range("J103").Value = 4375
$ws.Range("K103").Value = 15000
$ws.Range("L103").Value = 13125
$ws.Range("M103").Value = -14414
$ws.Range("N103").Value = -14297

# Row 113 (ALC)
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()

# Row 132 (ALC)
$ws.Range("H132").Value = 1705.1875
$ws.Range("I132").Value = 1377.3572
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 4132.071599999999
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -1602.071599999999
$ws.Range("N132").Value = -17060

# Row 138 (ALC)
$ws.Range("H138").Value = 4188.6665
$ws.Range("I138").Value = 3874.75
$ws.Range("J138").Value = 4251.45
$ws.Range("K138").Value = 11624.25
$ws.Range("L138").Value = 12754.35
$ws.Range("M138").Value = -6484.25
$ws.Range("N138").Value = -23034.35

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (ARM)
$ws.Range("H32").Value = 3516.2144
$ws.Range("I32").Value = 3167.025
$ws.Range("J32").Value = 10500
$ws.Range("K32").Value = 3167.025
$ws.Range("L32").Value = 10500
$ws.Range("M32").Value = -2880.025
$ws.Range("N32").Value = -11074

# Row 45 (ARM)
$ws.Range("H45").Value = 3500
$ws.Range("I45").Value = 2000
$ws.Range("J45").Value = 4000
$ws.Range("K45").Value = 2000
$ws.Range("L45").Value = 4000
$ws.Range("M45").Value = -1623
$ws.Range("N45").Value = -4754

# Row 74 (ARM)
$ws.Range("H74").Value = 837.46155
$ws.Range("I74").Value = 862.4545000000001
$ws.Range("J74").Value = 700
$ws.Range("K74").Value = 862.4545000000001
$ws.Range("L74").Value = 700
$ws.Range("M74").Value = 11.54549999999995
$ws.Range("N74").Value = -2448

# Row 77 (ARM)
$ws.Range("H77").Value = 837.46155
$ws.Range("I77").Value = 862.4545000000001
$ws.Range("J77").Value = 700
$ws.Range("K77").Value = 4312.2725
$ws.Range("L77").Value = 3500
$ws.Range("M77").Value = 55.72749999999996
$ws.Range("N77").Value = -12236

# Row 132 (ARM)
$ws.Range("H132").Value = 2213.4783
$ws.Range("I132").Value = 1060.8
$ws.Range("J132").Value = 4374.75
$ws.Range("K132").Value = 3182.4
$ws.Range("L132").Value = 13124.25
$ws.Range("M132").Value = -652.3999999999996
$ws.Range("N132").Value = -18184.25

$ws = $wb.Worksheets.Item("BSM")
# Row 107 (BSM)
$ws.Range("H107").Value = 899
$ws.Range("I107").Value = 798

# Row 116 (BSM)
$ws.Range("H116").Value = 40000
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 40000
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 40000
$ws.Range("N116").Value = -49178

# Row 134 (BSM)
$ws.Range("H134").Value = 10059
$ws.Range("I134").Value = 3412.2856
$ws.Range("J134").Value = 15874.875
$ws.Range("K134").Value = 10236.8568
$ws.Range("L134").Value = 47624.625
$ws.Range("M134").Value = -7701.856800000001
$ws.Range("N134").Value = -52694.625

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (CRP)
$ws.Range("H31").Value = 3280.8572
$ws.Range("I31").Value = 1244.3334
$ws.Range("J31").Value = 4808.25
$ws.Range("K31").Value = 1244.3334
$ws.Range("L31").Value = 4808.25
$ws.Range("M31").Value = -949.3334
$ws.Range("N31").Value = -5398.25

# Row 34 (CRP)
$ws.Range("H34").Value = 3280.8572
$ws.Range("I34").Value = 1244.3334
$ws.Range("J34").Value = 4808.25
$ws.Range("K34").Value = 1244.3334
$ws.Range("L34").Value = 4808.25
$ws.Range("M34").Value = -1042.3334
$ws.Range("N34").Value = -5212.25

# Row 86 (CRP)
$ws.Range("H86").Value = 8867.817999999999
$ws.Range("I86").Value = 5771.8
$ws.Range("J86").Value = 11447.833
$ws.Range("K86").Value = 5771.8
$ws.Range("L86").Value = 11447.833
$ws.Range("M86").Value = -4648.8
$ws.Range("N86").Value = -13693.833

# Row 89 (CRP)
$ws.Range("H89").Value = 8867.817999999999
$ws.Range("I89").Value = 5771.8
$ws.Range("J89").Value = 11447.833
$ws.Range("K89").Value = 28859
$ws.Range("L89").Value = 57239.165
$ws.Range("M89").Value = -23243
$ws.Range("N89").Value = -68471.16500000001

# Row 134 (CRP)
$ws.Range("H134").Value = 5000
$ws.Range("I134").Value = 5000
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 15000
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -12465

$ws = $wb.Worksheets.Item("CUL")
# Row 8 (CUL)
$ws.Range("H8").Value = 1001
$ws.Range("I8").Value = 1001
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 3003
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -2864

# Row 107 (CUL)
$ws.Range("H107").Value = 739.8
$ws.Range("I107").Value = 900
$ws.Range("J107").Value = 699.75
$ws.Range("K107").Value = 2700
$ws.Range("L107").Value = 2099.25
$ws.Range("M107").Value = -780
$ws.Range("N107").Value = -5939.25

$ws = $wb.Worksheets.Item("GSM")
# Row 102 (GSM)
$ws.Range("H102").Value = 2170
$ws.Range("I102").Value = 2170
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2170
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -548

# Row 126 (GSM)
$ws.Range("H126").Value = 3365.2
$ws.Range("I126").Value = 2470.6667
$ws.Range("J126").Value = 4707
$ws.Range("K126").Value = 7412.000100000001
$ws.Range("L126").Value = 14121
$ws.Range("M126").Value = -4942.000100000001
$ws.Range("N126").Value = -19061

# Row 132 (GSM)
$ws.Range("H132").Value = 3734.5
$ws.Range("I132").Value = 3102.25
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 9306.75
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -6776.75
$ws.Range("N132").Value = -20057

$ws = $wb.Worksheets.Item("LTW")
# Row 16 (LTW)
$ws.Range("H16").Value = 1156.5
$ws.Range("I16").Value = 987.8
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 987.8
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -817.8
$ws.Range("N16").Value = -2340

# Row 100 (LTW)
$ws.Range("H100").Value = 1016.4545
$ws.Range("I100").Value = 1079.1
$ws.Range("J100").Value = 390
$ws.Range("K100").Value = 1079.1
$ws.Range("L100").Value = 390
$ws.Range("M100").Value = -538.0999999999999
$ws.Range("N100").Value = -1472

# Row 132 (LTW)
$ws.Range("H132").Value = 3569.1538
$ws.Range("I132").Value = 3323
$ws.Range("J132").Value = 4123
$ws.Range("K132").Value = 9969
$ws.Range("L132").Value = 12369
$ws.Range("M132").Value = -7439
$ws.Range("N132").Value = -17429

$ws = $wb.Worksheets.Item("WVR")
# Row 107 (WVR)
$ws.Range("H107").Value = 531.1667
$ws.Range("I107").Value = 531.1667
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1593.5001
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 326.4999

# Row 132 (WVR)
$ws.Range("H132").Value = 3527.6216
$ws.Range("I132").Value = 2230.5833
$ws.Range("J132").Value = 5922.154
$ws.Range("K132").Value = 6691.749899999999
$ws.Range("L132").Value = 17766.462
$ws.Range("M132").Value = -4161.749899999999
$ws.Range("N132").Value = -22826.462

# Row 136 (WVR)
$ws.Range("H136").Value = 3500
$ws.Range("I136").Value = 3750
$ws.Range("J136").Value = 3250
$ws.Range("K136").Value = 11250
$ws.Range("L136").Value = 9750
$ws.Range("M136").Value = -8700
$ws.Range("N136").Value = -14850
